# Update gh-pages generated output figures (想去人数 / 最低票价) across sheets.
$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5437
$ws1.Range("F5").Value = 307
$ws1.Range("F7").Value = 19
$ws1.Range("F8").Value = 347

# Sheet: 演出 (Show)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 43
$ws2.Range("G3").Value = 80

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5437
$ws4.Range("F5").Value = 307
$ws4.Range("F7").Value = 19
$ws4.Range("F8").Value = 43
$ws4.Range("F9").Value = 347
$ws4.Range("G11").Value = 80
